# Apply updated dSF (column F) values to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    9  = -3
    13 = -3
    14 = 1
    16 = -8
    17 = -7
    18 = -2
    20 = -4
    23 = 0
    27 = 3
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
